$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of work log data
$ws.Range("A22").Value = 44413
$ws.Range("B22").Value = 6
$ws.Range("D22").Value = "Improved data annotations. Started to improve navigation of charts"

$ws.Range("A23").Value = 44428
$ws.Range("B23").Value = 7
$ws.Range("D23").Value = "Research on canvas and experimenting mit UI/UX design."

# Update the active selection
$ws.Range("H16").Select()
